$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells originally hold text (inline string) values, including
# cells in column D that look like plain numbers (e.g. "1.007", "308.26").
# Force those to stay text so Excel does not silently coerce them to
# floating point numbers, then restore the default "Normal" style so no
# stray number-format style gets introduced.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.509.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.810.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.007"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4559"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3663"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07125"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8765"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07777"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.821.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.283"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.362"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.37%  "
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008608"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.597.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.986"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.063"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "112.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.856"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08682"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.532"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7338"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.117"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.640"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.92%  "
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.084"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01946"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.908"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05114"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.994"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5018"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1558"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.167"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.007"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4617"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.982"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.589"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.70%  "
